$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.281.44"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "'1.911.31"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'327.78"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.4620"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").Value = "'0.3951"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "'0.07924"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "'22.38"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").Value = "'1.970.86"
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("D13").Value = "'7.091"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'5.765"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "'0.06970"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'88.69"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'0.00001006"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'17.13"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'29.276.32"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'5.345"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").Value = "'11.09"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'2.234.95"
$ws.Range("E24").Value = "  +6.41%  "
$ws.Range("D25").Value = "'2.068"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "'156.87"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").Value = "'19.45"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "'6.126"
$ws.Range("E28").Value = "  +7.44%  "
$ws.Range("D29").Value = "'1.977"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").Value = "'118.39"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'0.09389"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "'0.9235"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D33").Value = "'5.351"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "'1.358"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "'3.276"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'1.198"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("D37").Value = "'0.05841"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "'0.02107"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "'7.924"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("D40").Value = "'1.003"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "'0.5751"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "'0.1800"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "'9.957"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "'2.301"
$ws.Range("E44").Value = "  +8.27%  "
$ws.Range("D45").Value = "'12.00"
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("D46").Value = "'0.5418"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("D47").Value = "'0.07060"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'1.877"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").Value = "'2.554"
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("D50").Value = "'113.46"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  -2.64%  "
